$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(222).Insert()

$ws.Cells.Item(222, 1).Value = 5
$ws.Cells.Item(222, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(222, 3).Value = "Maule"
$ws.Cells.Item(222, 4).Value = 44505
$ws.Cells.Item(222, 5).Value = 7
$ws.Cells.Item(222, 6).Value = 100112043
$ws.Cells.Item(222, 7).Value = "Pepino ensalada"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 500
$ws.Cells.Item(222, 11).Value = 7000
$ws.Cells.Item(222, 12).Value = 7000
$ws.Cells.Item(222, 13).Value = 7000
$ws.Cells.Item(222, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(222, 15).Value = "Región del Maule"
$ws.Cells.Item(222, 16).Value = 88
$ws.Cells.Item(222, 17).Value = 80
$ws.Cells.Item(222, 18).Value = "Hortaliza"
